$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "네이버 블로그 원래 구글 검색에 안 뜸?"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/223217882577"

$ws.Range("D32").Value = "Feature Selection :: Recursive Feature Elimination (RFE)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/462"

$ws.Range("D51").Value = "[opensearch] text 필드에 하이픈이 들어간 경우 term 서치가 안 되는 이슈 해결"
$ws.Range("E51").Value = "https://bskyvision.com/entry/opensearch-text-%ED%95%84%EB%93%9C%EC%97%90-%ED%95%98%EC%9D%B4%ED%94%84%EC%9D%B4-%EB%93%A4%EC%96%B4%EA%B0%84-%EA%B2%BD%EC%9A%B0-term-%EC%84%9C%EC%B9%98%EA%B0%80-%EC%95%88-%EB%90%98%EB%8A%94-%EC%9D%B4%EC%8A%88-%ED%95%B4%EA%B2%B0"
